$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unraveling the Human Genome's Secrets" "The Profound Impact of Scientific Exploration: Unlocking the Wonders of Our Universe"

# --- Author name: " Sarah J" + "." + " Wilson" -> " Sam Taylor" ---
Replace-Text " Sarah J. Wilson" " Sam Taylor"

# --- Email: "genetics@advancescience" -> "sam" ; "org" -> "taylor@highschool" + "." + "edu" ---
Replace-Text "genetics@advancescience" "sam"
Replace-Text "org" "taylor@highschool.edu"

# --- Body paragraph (essay) ---
Replace-Text "With the advent of the Human Genome Project, the world was thrust into an exhilarating era of scientific discovery" "The realm of science beckons us with its enigmatic allure, inviting us to embark on a quest to unravel the tapestry of mysteries that enshroud our world"

Replace-Text " A collaboration of scientists, fueled by insatiable curiosity and a fervent desire to understand the intricacies of human existence, embarked on a colossal expedition to decipher the genetic blueprint of our species" " It is the avenue we have chosen to understand our place within this vast universe"

Replace-Text " As the vast ocean of genetic information yielded to their dedicated efforts, remarkable insights into the fundamental mechanisms of life unveiled themselves" " Through scientific exploration, we delve into the intricate mechanisms of nature, unravelling the secrets that lie hidden beneath its surface, and weaving together a symphony of knowledge"

Replace-Text "This momentous undertaking revealed a mesmerizingly intricate dance of three billion pairs of nucleotides gracefully twirling within the nucleus of every cell, holding the codes of life" "With each experiment conducted and each observation meticulously recorded, we step closer towards comprehending the profound interconnectedness of all things"

Replace-Text " As researchers delved deeper into this awe-inspiring tapestry of DNA, they uncovered genes, the fundamental units of heredity, each possessing a specific set of instructions for building and sustaining our bodies" " From the delicate dance of atoms to the grand cosmic symphony of stars, science sheds light upon the intricate workings of the universe, revealing a tapestry of interrelationships that bind us inextricably to our surroundings"

Replace-Text " This profound revelation laid bare the intricate interplay between genes and environment, shaping the symphony of our existence from the genetic score we inherit" " As we journey through the intricate pathways of scientific exploration, we unveil the hidden wonders that shape our existence, gaining a deeper appreciation for the inherent beauty and boundless mysteries that permeate our world"

Replace-Text "Furthermore, the Human Genome Project illuminated the powerful genetic basis of health and disease, revolutionizing the field of medicine" "Furthermore, scientific exploration fuels innovation, propelling society forward with advancements that transform our lives"

Replace-Text " By unraveling the genetic underpinnings of ailments, scientists gained the ability to diagnose, treat, and even prevent a myriad of conditions with remarkable precision" " The pursuit of knowledge transcends mere intellectual curiosity; it holds the potential to address global challenges, improve human well-being, and shape a future where progress and prosperity prevail"

# Merge of: " This transformative knowledge..." + "." + " With each discovery..." -> single sentence
Replace-Text " This transformative knowledge propelled the development of gene therapies, enabling the prospect of targeting and correcting genetic defects at their source. With each discovery, the tapestry of human biology grew richer, orchestrating a mesmerizing symphony of understanding that forever changed our perception of ourselves and our place in the grand symphony of life" " In this realm of exploration, we find hope, inspiration, and limitless possibilities for a world where knowledge reigns supreme"

# --- Summary paragraph ---
Replace-Text "The Human Genome Project, a monumental endeavor in scientific history, unveiled the astonishing secrets of our genetic heritage" "Dr"

Replace-Text " It illuminated the intricate interplay between genes and the environment, orchestrating the symphony of our existence" " Sam Taylor's essay, ""The Profound Impact of Scientific Exploration: Unlocking the Wonders of Our Universe,"" elucidates the transformative power of scientific exploration in unveiling the mysteries of our universe"

Replace-Text " This profound knowledge revolutionized medicine, transforming our ability to diagnose, treat, and prevent diseases with unprecedented precision" " The essay emphasizes the inherent beauty of the scientific journey, where experimentation and observation lead us towards a deeper comprehension of nature's intricate interrelationships"

Replace-Text " The tapestry of human biology, once shrouded in mystery, now " " It underscores the significance of scientific exploration in fueling innovation, addressing global challenges, and shaping a future driven by progress and well-"

Replace-Text "unfolds before us, revealing the profound influence of genetics on our health, heritage, and potential" "being"

Replace-Text " The seeds sown by the Human Genome Project continue to bear fruit, propelling us towards a future where genetic insights empower us to improve human health and well-being" " Through the exploration of science, we embark on an extraordinary adventure, unearthing the wonders of our universe and enriching our lives with knowledge and understanding"

# --- Add a new empty paragraph at the end of the document (after the Summary paragraph, before sectPr) ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
